# Add a new "leaderboard and achievements" time-log entry to Sheet1.
#
# The commit adds one new row (row 56) below the existing log on Sheet1:
#   A56 = "2012.7.28+29"                         (date/range column)
#   B56 = "增加game center的排行榜和成就"          (goal/what-was-done column)
#   D56 = 6                                       (hours worked)
#
# and updates the current selection in the sheet view from D58 to C63.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a brand-new row at position 56. Using Insert() (rather than just
# writing into previously-empty cells) makes Excel carry the formatting
# of the row above (row 55, style index "1") down onto the new row, just
# like it would if a user selected row 55's format painter / inserted a
# row from the Excel UI.
[void]$ws.Rows.Item(56).Insert()

# Fill in the three populated cells of the new row.
$ws.Cells.Item(56, 1).Value = "2012.7.28+29"
$ws.Cells.Item(56, 2).Value = "增加game center的排行榜和成就"
$ws.Cells.Item(56, 4).Value = 6

# Keep the sheet scrolled to where it was (top-left around row 37) and
# move the active selection to C63, matching the saved view state.
$win = $excel.ActiveWindow
$win.ScrollRow = 37
$win.ScrollColumn = 1
[void]$ws.Range("C63").Select()
